$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.206803679466248
$ws.Range("B1").Value = 2.256045341491699
$ws.Range("C1").Value = 3.509639501571655
$ws.Range("D1").Value = 2.351178169250488
$ws.Range("E1").Value = 1.286410570144653
